$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Center-align the header/legend cells (B2:C4 block and the D6:K6 legend row)
$ws.Range("B2:C4").HorizontalAlignment = -4108
$ws.Range("D6:K6").HorizontalAlignment = -4108

# Narrow column C now that it no longer needs to fit the old wider label
$ws.Columns("C").ColumnWidth = 3.4

# New practice row (class AM 9:00 - 12:50): an extra 8-bit pattern example
$ws.Range("D25").Value = 1
$ws.Range("E25").Value = 1
$ws.Range("F25").Value = 0
$ws.Range("G25").Value = 0
$ws.Range("H25").Value = 1
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("N25").Value = 200

# Move the active selection down to the next empty answer row
$ws.Range("D27:H27").Select()
